$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the split "Acceptatiecriteria" / ":" runs into one run.
#    This happens in several ListParagraph bullets. Some are already
#    a single run (left untouched); four of them still have the text
#    split across two runs and need merging, one of which also has a
#    stray <w:lastRenderedPageBreak/> that should disappear.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -eq "Acceptatiecriteria:`r") {
        $rng = $para.Range
        $rng.Find.Execute("Acceptatiecriteria:", $false, $false, $false, $false, $false, $true, 1, $false, "Acceptatiecriteria:", 2) | Out-Null
    }
    elseif ($text -eq "Acceptatiecriteria: `r") {
        $rng = $para.Range
        $rng.Find.Execute("Acceptatiecriteria: ", $false, $false, $false, $false, $false, $true, 1, $false, "Acceptatiecriteria: ", 2) | Out-Null
    }
}

# ------------------------------------------------------------------
# 2) Rewrite the "actief(...)" user story sentence into the new
#    "meer dan 5 likes heb gegeven" wording, keeping the leading
#    run intact and replacing everything after it with freshly
#    constructed runs (including the spellStart/spellEnd proofErr
#    pair around "likes").
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*mee rechten krijg als ik actief*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range
    $replaceRange = $d.Range($full.Start, $full.End - 1)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
        '<w:p>' + `
        '<w:r><w:t xml:space="preserve">Als gebruiker wil ik mee rechten krijg als ik </w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">meer dan 5 </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>likes</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> heb gegeven </w:t></w:r>' + `
        '</w:p>' + `
        '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $replaceRange.InsertXML($xml)
}
